$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.421.35'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.846.90'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.00'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4705'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('D8').Value = '0.2743'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').Value = '0.06319'
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = '18.35'
$ws.Range('E10').Value = '  +14.17%  '
$ws.Range('D11').Value = '1.834.70'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '0.07447'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '4.935'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '84.92'
$ws.Range('E14').Value = '  +1.81%  '
$ws.Range('D15').Value = '0.6217'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '262.96'
$ws.Range('E16').Value = '  +15.46%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '30.362.65'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '12.67'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007322'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '0.9989'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '4.921'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = '5.875'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '165.02'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '8.927'
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').Value = '18.35'
$ws.Range('E26').Value = '  +3.72%  '
$ws.Range('D27').Value = '1.884'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').Value = '0.1023'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').Value = '1.349'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').Value = '4.014'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').Value = '3.825'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').Value = '0.04799'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('D33').Value = '1.131'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('D34').Value = '0.6923'
$ws.Range('E34').Value = '  -2.82%  '
$ws.Range('D35').Value = '2.685'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').Value = '0.01863'
$ws.Range('E36').Value = '  +2.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.670'
$ws.Range('E37').Value = '  +1.64%  '
$ws.Range('D38').Value = '0.8728'
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').Value = '1.962'
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').Value = '105.79'
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '0.4044'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').Value = '5.431'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.120'
$ws.Range('E44').Value = '  +2.51%  '
$ws.Range('D45').Value = '61.98'
$ws.Range('E45').Value = '  +4.16%  '
$ws.Range('D46').Value = '0.1193'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').Value = '33.76'
$ws.Range('E47').Value = '  +4.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.530'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05490'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '1.339'
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('D51').Value = '0.3664'
$ws.Range('E51').Value = '  +1.17%  '
